# Fix crash when enabling JudgmentVisuals for the first time
#
# The JudgmentVisuals translation sheet was missing several Spanish (and one
# Korean) translation strings, which caused a crash the first time the
# tweak was enabled because the generator could not find a translation for
# those keys. This adds the missing translations:
#   - ERROR_METER_TICK_LIFE (Spanish)
#   - ERROR_METER_TICK_SECONDS (Spanish)
#   - ERROR_METER_SENSITIVITY (Spanish)
#   - HIDE_PERFECTS (Korean + Spanish)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JudgmentVisuals")

# ERROR_METER_TICK_LIFE row -> add missing Spanish translation
$ws.Range("D6").Value = "Duración de las marcas en pantalla:"

# ERROR_METER_TICK_SECONDS row -> add missing Spanish translation
$ws.Range("D7").Value = "{0} segundo(s)"

# ERROR_METER_SENSITIVITY row -> add missing Spanish translation
$ws.Range("D8").Value = "Sensibilidad de la flecha:"

# HIDE_PERFECTS row -> add missing Korean and Spanish translations
$ws.Range("C11").Value = "`"정확`" 판정 숨기기"
$ws.Range("D11").Value = "Esconder juicios de `"!Perfecto!`""
